$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.702.10'
$ws.Range("E2").Value = '  +0.65%  '
$ws.Range("D3").Value = '2.708.33'
$ws.Range("E3").Value = '  +2.28%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '599.70'
$ws.Range("E5").Value = '  +0.39%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '161.69'
$ws.Range("E6").Value = '  +3.17%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +0.15%  '
$ws.Range("D9").Value = '2.707.72'
$ws.Range("E9").Value = '  +2.29%  '
$ws.Range("E10").Value = '  +0.47%  '
$ws.Range("E12").Value = '  +1.01%  '
$ws.Range("E13").Value = '  +3.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.44'
$ws.Range("E14").Value = '  +1.44%  '
$ws.Range("D15").Value = '3.194.02'
$ws.Range("E15").Value = '  +2.05%  '
$ws.Range("E16").Value = '  -0.69%  '
$ws.Range("D17").Value = '68.616.78'
$ws.Range("E17").Value = '  +0.55%  '
$ws.Range("D18").Value = '2.708.96'
$ws.Range("E18").Value = '  +2.43%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.87'
$ws.Range("E19").Value = '  +4.44%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.67'
$ws.Range("E20").Value = '  +4.53%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '365.23'
$ws.Range("E21").Value = '  +0.48%  '
$ws.Range("E22").Value = '  +3.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.92'
$ws.Range("E23").Value = '  +2.44%  '
$ws.Range("E24").Value = '  +2.55%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '74.18'
$ws.Range("E25").Value = '  -1.34%  '
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.95'
$ws.Range("E27").Value = '  +2.01%  '
$ws.Range("E28").Value = '  +2.11%  '
$ws.Range("E29").Value = '  +1.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '598.06'
$ws.Range("E30").Value = '  +6.82%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  -0.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.29'
$ws.Range("E32").Value = '  +2.86%  '
$ws.Range("E33").Value = '  +2.71%  '
$ws.Range("E34").Value = '  +4.73%  '
$ws.Range("E35").Value = '  +3.29%  '
$ws.Range("E36").Value = '  +5.50%  '
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '161.49'
$ws.Range("E38").Value = '  -0.24%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.87'
$ws.Range("E39").Value = '  +0.91%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.381'
$ws.Range("E40").Value = '  +2.40%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.92'
$ws.Range("E41").Value = '  +2.37%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.43'
$ws.Range("E42").Value = '  +1.80%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.69'
$ws.Range("E43").Value = '  +3.32%  '
$ws.Range("D45").Value = '0.0₆0318'
$ws.Range("E45").Value = '  -5.35%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '158.17'
$ws.Range("E47").Value = '  -0.51%  '
$ws.Range("E48").Value = '  +5.83%  '
$ws.Range("E49").Value = '  +5.45%  '
$ws.Range("E50").Value = '  +7.28%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '22.14'
$ws.Range("E51").Value = '  +0.16%  '
